$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "Bronze rings" component row (row 8) to the FR_A0600 clutch
# actuation system BOM.
$ws.Range("C8").Value = "Bronze rings"
$ws.Range("D8").Value = "b"
$ws.Range("E8").Value = "Lubrificating rings for the lever"
$ws.Range("F8").Value = 2
$ws.Range("G8").Value = "FR_06006"

# Scroll the view over and move the active selection onto the newly added
# row, mirroring the author's on-screen focus when they made the edit.
$excel.ActiveWindow.ScrollColumn = 3
$ws.Range("G8").Select()
